$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 118
$ws1.Range("F5").Value = 20
$ws1.Range("F6").Value = 490
$ws1.Range("F8").Value = 104
$ws1.Range("F10").Value = 6461
$ws1.Range("F13").Value = 2613
$ws1.Range("F14").Value = 157
$ws1.Range("F15").Value = 267
$ws1.Range("F17").Value = 499

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 118
$ws4.Range("F7").Value = 20
$ws4.Range("F8").Value = 490
$ws4.Range("F10").Value = 104
$ws4.Range("F13").Value = 6461
$ws4.Range("F17").Value = 2613
$ws4.Range("F18").Value = 157
$ws4.Range("F19").Value = 267
$ws4.Range("F21").Value = 499
